$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.034.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = '''1.637.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.54%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''214.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.23%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = '''  -0.04%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.53%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -1.58%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  -1.42%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''18.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -4.15%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.0794'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +0.32%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''1.713.21'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +3.68%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''4.21'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -1.46%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = '''  -1.87%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.0₃0750'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -1.90%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''62.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.80%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''26.059.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +0.42%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  +0.52%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''191.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -0.71%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''  -1.65%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''  -2.79%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  -1.55%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  +1.56%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''143.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +0.49%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = '''Toncoin'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = '''1.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -0.72%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = '''BinanceUSD'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = '''https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = '''1.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +0.59%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''6.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.45%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''15.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -1.66%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E30").Value = '''  -2.85%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  -2.24%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -2.91%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  -1.55%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  -0.78%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.880'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -2.19%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''1.131.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.18%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  +0.09%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.528'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -2.49%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '''  -0.80%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''99.03'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.22%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.788'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -1.11%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''5.32'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -2.94%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -0.77%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''55.61'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -1.64%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -0.50%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  +1.65%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = '''7.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -1.20%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  +0.29%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.0931'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -2.76%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  +0.00%  '
$ws.Range("E51").Style = "Normal"
